$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.379.22"
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.227.35"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.82"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.78"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.403"
$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0906"
$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.552.88"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.55"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.24"
$ws.Range("E14").Value = "  +3.48%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.60"
$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.798"
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.235.33"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.205.79"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0945"
$ws.Range("E19").Value = "  +6.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  +2.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.13"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "244.44"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("E24").Value = "  +3.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +1.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.71"
$ws.Range("E26").Value = "  +1.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.60"
$ws.Range("E27").Value = "  +1.36%  "

$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.40"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +4.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  -0.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("E33").Value = "  +0.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.63"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0652"
$ws.Range("E35").Value = "  +5.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.36"
$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  -3.50%  "

$ws.Range("E39").Value = "  +6.44%  "

$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000230"
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.56"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0961"
$ws.Range("E43").Value = "  -1.46%  "

$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.08"
$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("E46").Value = "  -8.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.455.78"
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "15.91"
$ws.Range("E50").Value = "  -3.30%  "

$ws.Range("E51").Value = "  +3.59%  "
